{"js": "// Preparing for a staging release.docx \u2014 remove the obsolete \"Inspect the\n// CloudSdkSyncSample directory...\" paragraph, and relocate the Word\n// \"_GoBack\" bookmark (last-edit-location marker) from in front of the first\n// \"Git branch - -d <branch>\" bullet to the start of the paragraph that now\n// immediately follows \"Preparing for a staging release:\" (\"Exit Visual\n// Studio in the Cloud_SDK solution.\").\n\nconst body = context.document.body;\n\n// 1) Locate and delete the \"Inspect the CloudSdkSyncSample ...\" paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet inspectParagraph = null;\nlet exitParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (inspectParagraph === null && text.indexOf(\"Inspect the\") === 0) {\n    inspectParagraph = paragraphs.items[i];\n  }\n  if (exitParagraph === null && text.indexOf(\"Exit Visual Studio\") === 0) {\n    exitParagraph = paragraphs.items[i];\n  }\n}\n\nif (inspectParagraph) {\n  inspectParagraph.delete();\n}\n\n// 2) Move the \"_GoBack\" bookmark: delete it from wherever it currently is,\n//    then re-insert it at the very start of the \"Exit Visual Studio...\"\n//    paragraph.\nconst doc = context.document;\ndoc.deleteBookmark(\"_GoBack\");\n\nif (exitParagraph) {\n  const startRange = exitParagraph.getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# \"Preparing for a staging release.docx\" maintenance edit:\n#   1) Remove the obsolete \"Inspect the CloudSdkSyncSample directory\n#      hierarchy ...\" paragraph entirely.\n#   2) Relocate the \"_GoBack\" bookmark (Word's last-edit-location marker)\n#      from in front of the first \"Git branch - -d <branch>\" bullet to the\n#      very start of the paragraph that now immediately follows\n#      \"Preparing for a staging release:\" -- i.e. \"Exit Visual Studio in\n#      the Cloud_SDK solution.\"\n\n$d = $word.ActiveDocument\n\n# 1) Delete the \"Inspect the ...\" paragraph.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Inspect the*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2) Move the \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Exit Visual Studio*\") {\n        $r = $p.Range.Duplicate\n        $r.Collapse(1)   # wdCollapseStart\n        $d.Bookmarks.Add(\"_GoBack\", $r)\n        break\n    }\n}\n"}
